$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New localization row: InsufficientFunds
$ws.Cells.Item(91, 1).Value = "InsufficientFunds"
$ws.Cells.Item(91, 2).Value = "SMS bakiyeniz yetersizdir."

# Match the formatting used by the header row (both columns share style index 1)
$ws.Range("A1:B1").Copy()
$ws.Range("A91:B91").PasteSpecial(-4122)
